$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: was a "Variable" row for MHC A3 (242...) -> becomes a "Constant" row for B2M ("dddd")
$ws.Range("B7").Value = "Constant"
$ws.Range("C7").Value = "B2M"
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = "dddd"
$ws.Range("J7").ClearContents()

# Row 8: was a "Constant" row for MHC A3 ("dddd") -> becomes a "Variable" row for B2M (242...)
$ws.Range("B8").Value = "Variable"
$ws.Range("C8").Value = "B2M"
$ws.Range("D8").Value = 242
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "Top"
$ws.Range("G8").Value = "GCAGGTGACGGTACT"
$ws.Range("H8").Value = "CAAAAGTGGGCAGCA"
$ws.Range("I8").Value = "W"
$ws.Range("J8").Value = "WT"

# Row 9: was a "Variable" row for MHC A3 (246...) -> becomes a "Constant" row for B2M ("eeee")
$ws.Range("B9").Value = "Constant"
$ws.Range("C9").Value = "B2M"
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = "eeee"
$ws.Range("J9").ClearContents()

# Row 10: was a "Constant" row for MHC A3 ("eeee") -> becomes a "Variable" row for B2M (246...)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Variable"
$ws.Range("C10").Value = "B2M"
$ws.Range("D10").Value = 246
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "Top"
$ws.Range("G10").Value = "GTTGTCGTTCCATCT"
$ws.Range("H10").Value = "CAAGAACAGCGTTAT"
$ws.Range("I10").Value = "G"
$ws.Range("J10").Value = "GT"

# Row 11: brand-new trailing "Constant" row for B2M ("ffff")
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Constant"
$ws.Range("C11").Value = "B2M"
$ws.Range("I11").Value = "ffff"

# Column B now has an explicit width, matching column I's existing width
$ws.Range("B1").ColumnWidth = 15.6666666666666

# Final selection, as recorded in the saved file
$ws.Range("G18").Select()
